$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Remove the old parameter table body (rows 8-40) and insert a fresh,
# larger block of 45 rows (8-52) to hold the restructured parameter list.
$ws.Rows("8:40").Delete()
$ws.Rows("8:52").Insert()

$ws.Range("B8").Value = "epi_proportion_cases"
$ws.Range("C8").Value = 1
$ws.Range("B9").Value = "tb_multiplier_force_smearpos"
$ws.Range("C9").Value = 1
$ws.Range("B10").Value = "tb_multiplier_force_smearneg"
$ws.Range("C10").Value = 0.24
$ws.Range("B11").Value = "tb_multiplier_force_extrapul"
$ws.Range("C11").Value = 0
$ws.Range("B12").Value = "tb_multiplier_force"
$ws.Range("C12").Value = 1
$ws.Range("B13").Value = "tb_n_contact"
$ws.Range("C13").Value = 140
$ws.Range("B14").Value = "tb_proportion_early_progression"
$ws.Range("C14").Value = 0.12
$ws.Range("B15").Value = "tb_timeperiod_early_latent"
$ws.Range("C15").Value = 0.4
$ws.Range("B16").Value = "tb_rate_late_progression"
$ws.Range("C16").Value = 0.007
$ws.Range("B17").Value = "tb_proportion_casefatality_untreated_smearpos"
$ws.Range("C17").Value = 0.7
$ws.Range("B18").Value = "tb_proportion_casefatality_untreated_smearneg"
$ws.Range("C18").Value = 0.2
$ws.Range("B19").Value = "tb_proportion_casefatality_untreated"
$ws.Range("C19").Value = 0.4
$ws.Range("B20").Value = "tb_timeperiod_activeuntreated"
$ws.Range("C20").Value = 3
$ws.Range("B21").Value = "tb_multiplier_bcg_protection"
$ws.Range("C21").Value = 0.5
$ws.Range("B22").Value = "program_prop_vac"
$ws.Range("C22").Value = 0.88
$ws.Range("B23").Value = "program_prop_unvac"
$ws.Range("C23").Formula = "=1-C22"
$ws.Range("B24").Value = "program_proportion_detect"
$ws.Range("C24").Value = 0.8
$ws.Range("B25").Value = "program_algorithm_sensitivity"
$ws.Range("C25").Value = 0.9
$ws.Range("B26").Value = "program_rate_start_treatment"
$ws.Range("C26").Value = 26
$ws.Range("B27").Value = "tb_timeperiod_treatment_ds"
$ws.Range("C27").Value = 0.5
$ws.Range("B28").Value = "tb_timeperiod_treatment_mdr"
$ws.Range("C28").Value = 2
$ws.Range("B29").Value = "tb_timeperiod_treatment_xdr"
$ws.Range("C29").Value = 3
$ws.Range("B30").Value = "tb_timeperiod_treatment_inappropriate"
$ws.Range("C30").Value = 3
$ws.Range("B31").Value = "tb_timeperiod_infect_ontreatment_ds"
$ws.Range("C31").Value = 0.035
$ws.Range("B32").Value = "tb_timeperiod_infect_ontreatment_mdr"
$ws.Range("C32").Formula = "=1/12"
$ws.Range("B33").Value = "tb_timeperiod_infect_ontreatment_xdr"
$ws.Range("C33").Formula = "=2/12"
$ws.Range("B34").Value = "tb_timeperiod_infect_ontreatment_inappropriate"
$ws.Range("C34").Value = 2
$ws.Range("B35").Value = "program_proportion_success_ds"
$ws.Range("C35").Value = 0.9
$ws.Range("B36").Value = "program_proportion_success_mdr"
$ws.Range("C36").Value = 0.6
$ws.Range("B37").Value = "program_proportion_success_xdr"
$ws.Range("C37").Value = 0.4
$ws.Range("B38").Value = "program_proportion_success_inappropriate"
$ws.Range("C38").Value = 0.25
$ws.Range("B39").Value = "program_rate_restart_presenting"
$ws.Range("C39").Value = 4
$ws.Range("B40").Value = "proportion_amplification"
$ws.Range("C40").Formula = "=1/15"
$ws.Range("B41").Value = "timepoint_introduce_mdr"
$ws.Range("C41").Value = 1950
$ws.Range("B42").Value = "timepoint_introduce_xdr"
$ws.Range("C42").Value = 2050
$ws.Range("B43").Value = "treatment_available_date"
$ws.Range("C43").Value = 1940
$ws.Range("B44").Value = "dots_start_date"
$ws.Range("C44").Value = 1990
$ws.Range("B45").Value = "finish_scaleup_date"
$ws.Range("C45").Value = 2010
$ws.Range("B46").Value = "pretreatment_available_proportion"
$ws.Range("C46").Value = 0
$ws.Range("B47").Value = "dots_start_proportion"
$ws.Range("C47").Value = 0.85
$ws.Range("B48").Value = "program_prop_assign_mdr"
$ws.Range("C48").Value = 0.6
$ws.Range("B49").Value = "program_prop_assign_xdr"
$ws.Range("C49").Value = 0.4
$ws.Range("B50").Value = "program_prop_lowquality"
$ws.Range("C50").Value = 0.4
$ws.Range("B51").Value = "program_rate_leavelowquality"
$ws.Range("C51").Value = 2
$ws.Range("B52").Value = "program_prop_nonsuccessoutcomes_death"
$ws.Range("C52").Value = 0.25

# Move the selection to match the authored state.
$ws.Range("C14").Select()
